# Generate Report for Handback
# Adds a new handback entry (dffb5ce4-5a26-49fb-b42e-b6af9bc09fba.md) as row 3
# on the Overview / zh-cn / de-de sheets, alongside updated timestamps/hashes
# for the existing 0755463a-7804-4d15-8686-2fe2c663ed24.md entry (renamed from
# 1b2ddbd0-62e0-448c-b24a-fc916ab68aaf.md).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldGuid = "1b2ddbd0-62e0-448c-b24a-fc916ab68aaf"
$guid1   = "0755463a-7804-4d15-8686-2fe2c663ed24"
$guid2   = "dffb5ce4-5a26-49fb-b42e-b6af9bc09fba"

$zhHash1 = "e25a99292201a7e940a41497072552aab17ed612"
$deHash1 = "e25a99292201a7e940a41497072552aab17ed612"
$zhHash2 = "d942a901efb7f92585b47c972aea5cb634d312dd"
$deHash2 = "d942a901efb7f92585b47c972aea5cb634d312dd"

# ---------------------------------------------------------------------------
# Step 1: rename the first GUID throughout (1b2ddbd0... -> 0755463a...)
# ---------------------------------------------------------------------------

$wsOverview.Range("A2").Value = "$guid1.md"
$wsOverview.Range("B2").Value = "e2e\$guid1.md"
$wsOverview.Range("G2").Value = "2016-08-30 15:20:41"

$wsZhCn.Range("A2").Value = "$guid1.md"
$wsZhCn.Range("I2").Value = "$guid1.md"
$wsZhCn.Range("G2").Value = "$guid1.$zhHash1.zh-cn.xlf"
$wsZhCn.Range("J2").Value = "$guid1.$zhHash1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-30 15:20:36"
$wsZhCn.Range("K2").Value = "2016-08-30 15:20:54"

$wsDeDe.Range("A2").Value = "$guid1.md"
$wsDeDe.Range("I2").Value = "$guid1.md"
$wsDeDe.Range("G2").Value = "$guid1.$deHash1.de-de.xlf"
$wsDeDe.Range("J2").Value = "$guid1.$deHash1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-30 15:20:41"
$wsDeDe.Range("K2").Value = "2016-08-30 15:21:03"

# ---------------------------------------------------------------------------
# Step 2: add the second GUID (dffb5ce4...) as a brand new row on every sheet
# ---------------------------------------------------------------------------

# -- Overview sheet (table3) --
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "$guid2.md"
$wsOverview.Range("B3").Value = "e2e\$guid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-30 15:20:41"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4543295d15b258bf1a0e069bce50f98f7b7ccf4c/e2e/$guid2.md", "", "", "e2e\$guid2.md") | Out-Null

# -- zh-cn sheet (table1) --
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = "$guid2.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = "$guid2.$zhHash2.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-30 15:20:36"
$wsZhCn.Range("I3").Value = "$guid2.md"
$wsZhCn.Range("J3").Value = "$guid2.$zhHash2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-30 15:20:54"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4543295d15b258bf1a0e069bce50f98f7b7ccf4c/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a28b48107aea743c27d9482920c8d0cf56d77671/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null

# -- de-de sheet (table2) --
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = "$guid2.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = "$guid2.$deHash2.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-30 15:20:41"
$wsDeDe.Range("I3").Value = "$guid2.md"
$wsDeDe.Range("J3").Value = "$guid2.$deHash2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-30 15:21:03"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4543295d15b258bf1a0e069bce50f98f7b7ccf4c/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/73caf943f634552e1329b402e93e8107499bc609/e2e/$guid2.md", "", "", "$guid2.md") | Out-Null

# Re-point the existing hyperlinks on row 2 (A2/I2/B2) to the renamed GUID so
# the display text matches the new file name everywhere.
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$guid1.md"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$guid1.md"
$wsZhCn.Hyperlinks.Item(2).TextToDisplay = "$guid1.md"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$guid1.md"
$wsDeDe.Hyperlinks.Item(2).TextToDisplay = "$guid1.md"
